$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "timestamp" column (O) for all data rows (2 through 26)
# from the old scrape time to the new scrape time.
for ($row = 2; $row -le 26; $row++) {
    $ws.Cells.Item($row, 15).Value = "2022-08-28 20:56:49"
}
